# SearchTestData.xlsx — "Adding dependency test cases for authoring /
# Changed validations for Type a head / Remove validations for search"
#
# The VALIDATIONS column (J) for the "search" dependency test rows used to
# hold long, composite validation strings (one clause per matched field,
# joined with "&&"). Those validations are replaced with the simple
# "status=200" check that is already used elsewhere in the sheet, and the
# rows' explicit (wrapped-text) row heights - which existed only to show the
# long strings - are cleared back to the sheet's default auto height.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose VALIDATIONS (column J) collapse down to the plain "status=200"
# check, and whose explicit row height goes away with them (autofit reverts
# to the sheet's default row height once the long text is gone).
$rowsToSimplify = @(2, 6, 7, 8, 9, 10, 11, 12, 13)

foreach ($r in $rowsToSimplify) {
    $ws.Range("J$r").Value = "status=200"
    $ws.Rows.Item($r).AutoFit()
}
